# Corrección en las complejidades para pasarlo a pdf
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update complexity values
$ws.Range("B9").Value = "O(n^2)"
$ws.Range("F10").Value = "O(2^(n))"

# Update the active cell selection
$ws.Range("H11").Select()

$wb.Save()
